$wb = $excel.ActiveWorkbook
$readme = $wb.Worksheets.Add()
$readme.Name = "README graph"
$readme.Move($wb.Worksheets.Item(1))

$readme.Range("B1").Value = "c405c58cbf (origin/master)"

$readme.Range("A2").Value = "Serialize::"
$labels = @("MessagePack for C#", "MessagePack for C# (LZ4)", "MsgPack-Cli", "protobuf-net", "ZeroFormatter", "Json.NET", "Json.NET(+GZip)")
$serialize = @(137.7, 160.73, 335.23, 449.62, 141.73, 1543.3, 1680.56)
$deserialize = @(177.33, 182.04, 1106.87, 627.23, 139.91, 1876.02, 2237.38)
$filesize = @(1803, 562, 2347, 2248, 5004, 6096, 458)

for ($i = 0; $i -lt 7; $i++) {
    $row = 3 + $i
    $readme.Range("A$row").Value = $labels[$i]
    $readme.Range("B$row").Value = $serialize[$i]
}

$readme.Range("A11").Value = "Deserialize::"
for ($i = 0; $i -lt 7; $i++) {
    $row = 12 + $i
    $readme.Range("A$row").Value = $labels[$i]
    $readme.Range("B$row").Value = $deserialize[$i]
}

$readme.Range("A20").Value = "FileSize::"
for ($i = 0; $i -lt 7; $i++) {
    $row = 21 + $i
    $readme.Range("A$row").Value = $labels[$i]
    $readme.Range("B$row").Value = $filesize[$i]
}

Write-Host "done"
